$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-28 Wednesday" "2026-01-29 Thursday"

Replace-Text "57÷7=" "61÷3="
Replace-Text "48÷3=" "25÷2="
Replace-Text "10÷6=" "23÷4="
Replace-Text "72÷8=" "65÷4="
Replace-Text "45÷8=" "94÷8="
Replace-Text "33÷4=" "39÷7="
Replace-Text "85÷3=" "78÷8="
Replace-Text "13÷6=" "56÷9="
Replace-Text "72÷6=" "68÷5="
Replace-Text "64÷5=" "53÷6="
Replace-Text "77÷8=" "79÷4="
Replace-Text "75÷5=" "29÷7="
Replace-Text "74÷2=" "59÷8="
Replace-Text "47÷2=" "72÷7="
Replace-Text "32÷7=" "89÷3="
Replace-Text "97÷8=" "24÷2="
Replace-Text "91÷2=" "27÷4="
Replace-Text "70÷3=" "47÷7="
Replace-Text "55÷7=" "65÷7="
Replace-Text "91÷7=" "79÷8="
Replace-Text "76÷6=" "96÷4="
Replace-Text "11÷6=" "81÷9="
Replace-Text "91÷6=" "14÷3="
Replace-Text "54÷3=" "96÷8="
Replace-Text "90÷3=" "82÷9="
